$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1980830670926517
$ws.Range("C2").Value = 0.5782747603833865
$ws.Range("J2").Value = 0.01277955271565495
$ws.Range("P2").Value = 0.1437699680511182
$ws.Range("S2").Value = 0.0670926517571885

# Row 3
$ws.Range("B3").Value = 0.01538461538461539
$ws.Range("C3").Value = 0.03589743589743589
$ws.Range("J3").Value = 0.02564102564102564
$ws.Range("P3").Value = 0.7128205128205128
$ws.Range("S3").Value = 0.2102564102564103

# Row 4
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.7441860465116279
$ws.Range("S4").Value = 0.2325581395348837

# Row 5
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333

# Row 6
$ws.Range("B6").Value = 0.04417670682730924
$ws.Range("D6").Value = 0.008032128514056224
$ws.Range("E6").Value = 0.004016064257028112
$ws.Range("F6").Value = 0.07630522088353414
$ws.Range("J6").Value = 0.2409638554216867
$ws.Range("Q6").Value = 0.144578313253012
$ws.Range("R6").Value = 0.05622489959839357
$ws.Range("S6").Value = 0.4257028112449799

# Row 7
$ws.Range("B7").Value = 0.06140350877192982
$ws.Range("D7").Value = 0.02192982456140351
$ws.Range("E7").Value = 0.004385964912280702
$ws.Range("F7").Value = 0.07456140350877193
$ws.Range("J7").Value = 0.1359649122807018
$ws.Range("O7").Value = 0.01754385964912281
$ws.Range("Q7").Value = 0.1052631578947368
$ws.Range("R7").Value = 0.07456140350877193
$ws.Range("S7").Value = 0.5043859649122807

# Row 8
$ws.Range("B8").Value = 0.09320388349514563
$ws.Range("D8").Value = 0.01941747572815534
$ws.Range("E8").Value = 0.001941747572815534
$ws.Range("F8").Value = 0.06407766990291262
$ws.Range("J8").Value = 0.1203883495145631
$ws.Range("O8").Value = 0.02912621359223301
$ws.Range("Q8").Value = 0.2097087378640777
$ws.Range("R8").Value = 0.09320388349514563
$ws.Range("S8").Value = 0.3689320388349515

# Row 9
$ws.Range("B9").Value = 0.1038062283737024
$ws.Range("D9").Value = 0.006920415224913495
$ws.Range("E9").Value = 0.01038062283737024
$ws.Range("F9").Value = 0.03114186851211072
$ws.Range("J9").Value = 0.1107266435986159
$ws.Range("O9").Value = 0.02422145328719723
$ws.Range("Q9").Value = 0.1660899653979239
$ws.Range("R9").Value = 0.08996539792387544
$ws.Range("S9").Value = 0.4567474048442907

# Row 10
$ws.Range("B10").Value = 0.09582477754962354
$ws.Range("D10").Value = 0.01779603011635866
$ws.Range("F10").Value = 0.06639288158795345
$ws.Range("J10").Value = 0.1286789869952088
$ws.Range("O10").Value = 0.01642710472279261
$ws.Range("Q10").Value = 0.1923340177960301
$ws.Range("R10").Value = 0.09308692676249145
$ws.Range("S10").Value = 0.3894592744695414

# Row 11
$ws.Range("G11").Value = 0.1232876712328767
$ws.Range("J11").Value = 0.07123287671232877
$ws.Range("K11").Value = 0.1780821917808219
$ws.Range("L11").Value = 0.6027397260273972
$ws.Range("S11").Value = 0.02465753424657534

# Row 12
$ws.Range("G12").Value = 0.7397260273972602
$ws.Range("J12").Value = 0.182648401826484
$ws.Range("L12").Value = 0.0045662100456621
$ws.Range("S12").Value = 0.0730593607305936

# Row 13
$ws.Range("G13").Value = 0.7105263157894737
$ws.Range("J13").Value = 0.1842105263157895
$ws.Range("S13").Value = 0.1052631578947368

# Row 15
$ws.Range("F15").Value = 0.01923076923076923
$ws.Range("H15").Value = 0.1769230769230769
$ws.Range("I15").Value = 0.08076923076923077
$ws.Range("J15").Value = 0.3192307692307692
$ws.Range("K15").Value = 0.05384615384615385
$ws.Range("M15").Value = 0.007692307692307693
$ws.Range("O15").Value = 0.06153846153846154
$ws.Range("S15").Value = 0.2807692307692308

# Row 16
$ws.Range("F16").Value = 0.01923076923076923
$ws.Range("H16").Value = 0.1778846153846154
$ws.Range("I16").Value = 0.1201923076923077
$ws.Range("J16").Value = 0.3990384615384616
$ws.Range("K16").Value = 0.09615384615384616
$ws.Range("M16").Value = 0.009615384615384616
$ws.Range("N16").Value = 0.004807692307692308
$ws.Range("O16").Value = 0.04807692307692308
$ws.Range("S16").Value = 0.125

# Row 17
$ws.Range("F17").Value = 0.01976284584980237
$ws.Range("H17").Value = 0.1600790513833992
$ws.Range("I17").Value = 0.1067193675889328
$ws.Range("J17").Value = 0.41699604743083
$ws.Range("K17").Value = 0.09881422924901186
$ws.Range("M17").Value = 0.0158102766798419
$ws.Range("N17").Value = 0.001976284584980237
$ws.Range("O17").Value = 0.06719367588932806
$ws.Range("S17").Value = 0.1126482213438735

# Row 18
$ws.Range("F18").Value = 0.01244813278008299
$ws.Range("H18").Value = 0.1701244813278008
$ws.Range("I18").Value = 0.1369294605809129
$ws.Range("J18").Value = 0.4481327800829876
$ws.Range("K18").Value = 0.07468879668049792
$ws.Range("M18").Value = 0.008298755186721992
$ws.Range("O18").Value = 0.04979253112033195
$ws.Range("S18").Value = 0.09958506224066389

# Row 19
$ws.Range("F19").Value = 0.01644528779253637
$ws.Range("H19").Value = 0.200506008855155
$ws.Range("I19").Value = 0.1005692599620493
$ws.Range("J19").Value = 0.3447185325743201
$ws.Range("K19").Value = 0.1214421252371917
$ws.Range("M19").Value = 0.01644528779253637
$ws.Range("N19").Value = 0.0006325110689437065
$ws.Range("O19").Value = 0.06451612903225806
$ws.Range("S19").Value = 0.1347248576850095
